# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E21) held the payroll-period labels for
# the debtor's overdue periods. The database refresh shifts the window of
# periods being reported: what used to read 2411/2410/2409/2408/2407 (desc.)
# now reads 2407/2408/2409/2410/2411 (asc.) while the last row (2412) and
# the middle row (2409) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2407"
$ws.Range("E17").Value = "2408"
$ws.Range("E19").Value = "2410"
$ws.Range("E20").Value = "2411"
